$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to remain a text value even if it looks numeric,
    # then restore the default (no explicit) style so no residual formatting is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '42.049.31'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '2.217.20'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("E4").Value = '  -0.02%  '

Set-TextValue $ws "D5" '242.95'
$ws.Range("E5").Value = '  -1.67%  '

$ws.Range("E6").Value = '  -0.46%  '

Set-TextValue $ws "D7" '73.56'
$ws.Range("E7").Value = '  -1.11%  '

$ws.Range("E8").Value = '  +0.17%  '

Set-TextValue $ws "D9" '0.613'
$ws.Range("E9").Value = '  -0.83%  '

Set-TextValue $ws "D10" '43.81'
$ws.Range("E10").Value = '  +6.15%  '

Set-TextValue $ws "D11" '0.0959'
$ws.Range("E11").Value = '  +2.00%  '

Set-TextValue $ws "D12" '7.11'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("D14").Value = '2.547.27'
$ws.Range("E14").Value = '  -1.40%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws "D15" '0.843'
$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws "D16" '14.22'
$ws.Range("E16").Value = '  -1.96%  '

$ws.Range("D17").Value = '2.261.36'
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '41.880.00'
$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("E19").Value = '  +12.56%  '

$ws.Range("E20").Value = '  +0.48%  '

Set-TextValue $ws "D21" '72.48'
$ws.Range("E21").Value = '  +0.83%  '

Set-TextValue $ws "D22" '10.50'
$ws.Range("E22").Value = '  +34.13%  '

$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("E24").Value = '  -8.05%  '

Set-TextValue $ws "D25" '11.57'
$ws.Range("E25").Value = '  +3.64%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  +1.29%  '

Set-TextValue $ws "D28" '2.28'
$ws.Range("E28").Value = '  -1.10%  '

Set-TextValue $ws "D29" '2.21'
$ws.Range("E29").Value = '  +6.04%  '

Set-TextValue $ws "D30" '166.64'
$ws.Range("E30").Value = '  -1.80%  '

Set-TextValue $ws "D31" '20.59'
$ws.Range("E31").Value = '  -0.26%  '

Set-TextValue $ws "D32" '5.62'
$ws.Range("E32").Value = '  +14.76%  '

Set-TextValue $ws "D33" '0.0796'
$ws.Range("E33").Value = '  -3.43%  '

$ws.Range("E34").Value = '  -0.26%  '

Set-TextValue $ws "D35" '29.38'
$ws.Range("E35").Value = '  -2.01%  '

$ws.Range("E36").Value = '  -4.17%  '

Set-TextValue $ws "D37" '4.29'
$ws.Range("E37").Value = '  -4.73%  '

$ws.Range("E38").Value = '  +0.36%  '

Set-TextValue $ws "D39" '12.95'
$ws.Range("E39").Value = '  -4.32%  '

$ws.Range("E40").Value = '  -2.45%  '

Set-TextValue $ws "D41" '64.82'
$ws.Range("E41").Value = '  +3.46%  '

Set-TextValue $ws "D42" '5.65'
$ws.Range("E42").Value = '  -2.35%  '

Set-TextValue $ws "D43" '0.199'
$ws.Range("E43").Value = '  -1.35%  '

Set-TextValue $ws "D44" '8.69'
$ws.Range("E44").Value = '  +0.32%  '

Set-TextValue $ws "D45" '104.03'
$ws.Range("E45").Value = '  -4.28%  '

$ws.Range("E46").Value = '  +0.74%  '

Set-TextValue $ws "D47" '2.39'
$ws.Range("E47").Value = '  +4.31%  '

$ws.Range("E48").Value = '  -0.75%  '

$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("E50").Value = '  +0.75%  '

$ws.Range("D51").Value = '2.424.11'
$ws.Range("E51").Value = '  -1.34%  '
